# Auto-generated Excel COM-interop script to update the cryptos price list
# (Price column D, Volume(1h) column E) for the "Updated cryptos list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sets a cell to a literal text value without letting Excel auto-convert
# numeric-looking text like "1.00" or "32.90" into a real number (which would
# silently drop the trailing zeros / original formatting). Note: positional
# parameters are used throughout this script instead of -Name Value syntax.
function Set-TextValue($Row, $Col, $Text) {
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# New "Price" values (column D = 4) keyed by row number.
$priceUpdates = [ordered]@{
    2 = "69.748.52"
    3 = "3.713.42"
    4 = "1.00"
    5 = "672.28"
    6 = "161.81"
    8 = "0.499"
    11 = "0.444"
    13 = "32.90"
    14 = "3.676.13"
    15 = "69.751.31"
    17 = "16.27"
    19 = "474.69"
    20 = "9.83"
    23 = "3.858.71"
    24 = "0.0000128"
    26 = "10.97"
    27 = "9.13"
    34 = "26.95"
    35 = "3.702.71"
    36 = "8.56"
    40 = "1.00"
    42 = "172.61"
    43 = "0.943"
    44 = "47.07"
    46 = "0.000282"
    47 = "27.77"
    50 = "7.90"
}

# New "Volume(1h)" values (column E = 5) keyed by row number.
# Values keep the original two leading / two trailing spaces.
$volumeUpdates = [ordered]@{
    2 = "  +0.52%  "
    3 = "  +0.95%  "
    4 = "  +0.06%  "
    5 = "  -1.91%  "
    6 = "  +1.80%  "
    7 = "  -0.03%  "
    8 = "  +0.90%  "
    9 = "  +0.50%  "
    10 = "  +0.63%  "
    11 = "  +1.99%  "
    12 = "  +1.16%  "
    14 = "  -0.36%  "
    15 = "  +0.57%  "
    16 = "  +1.52%  "
    17 = "  +2.75%  "
    18 = "  +1.88%  "
    19 = "  +0.99%  "
    20 = "  -1.39%  "
    21 = "  +0.89%  "
    22 = "  +0.86%  "
    23 = "  +0.89%  "
    24 = "  +4.00%  "
    25 = "  +0.00%  "
    26 = "  +0.18%  "
    27 = "  -0.40%  "
    28 = "  -0.53%  "
    29 = "  +0.38%  "
    30 = "  +1.69%  "
    31 = "  +1.07%  "
    32 = "  +4.79%  "
    33 = "  +1.17%  "
    34 = "  +0.21%  "
    35 = "  +1.32%  "
    36 = "  +4.44%  "
    37 = "  -0.57%  "
    38 = "  +0.01%  "
    39 = "  +1.41%  "
    40 = "  +0.00%  "
    41 = "  +1.28%  "
    42 = "  +3.90%  "
    43 = "  +0.28%  "
    45 = "  +2.86%  "
    46 = "  -0.29%  "
    47 = "  +2.08%  "
    48 = "  -0.87%  "
    49 = "  -1.52%  "
    50 = "  +1.39%  "
    51 = "  +1.12%  "
}

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $row 4 $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

